$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: issue number and week-covering dates ---
$issueCell = $ws.Range("A8")
$issueText = $issueCell.Value2
$issueLen = $issueText.Length
$issueCell.Characters($issueLen - 1, 2).Text = "20"

$dateCell = $ws.Range("C9")
$dateText = $dateCell.Value2
$dateText = $dateText.Replace("5/8/2023", "5/15/2023").Replace("5/14/2023", "5/21/2023")
$dateCell.Value2 = $dateText

# --- Weekly crime statistics table (rows 14-30) ---
# Row 14
$ws.Range("C14").Value2 = 4
$ws.Range("D14").Value2 = 5
$ws.Range("E14").Value2 = -20
$ws.Range("G14").Value2 = 14
$ws.Range("H14").Value2 = -35.714285714285
$ws.Range("I14").Value2 = 45
$ws.Range("J14").Value2 = 59
$ws.Range("K14").Value2 = -23.728813559322
$ws.Range("L14").Value2 = -16.666666666666
$ws.Range("M14").Value2 = 12.5
$ws.Range("N14").Value2 = -75.138121546961
# Row 15
$ws.Range("C15").Value2 = 4
$ws.Range("D15").Value2 = 6
$ws.Range("E15").Value2 = -33.333333333333
$ws.Range("F15").Value2 = 36
$ws.Range("H15").Value2 = 38.461538461538
$ws.Range("I15").Value2 = 159
$ws.Range("J15").Value2 = 150
$ws.Range("K15").Value2 = 6
$ws.Range("L15").Value2 = 22.307692307692
$ws.Range("M15").Value2 = 45.871559633027
$ws.Range("N15").Value2 = -40
# Row 16
$ws.Range("C16").Value2 = 103
$ws.Range("D16").Value2 = 85
$ws.Range("E16").Value2 = 21.176470588235
$ws.Range("F16").Value2 = 347
$ws.Range("G16").Value2 = 370
$ws.Range("H16").Value2 = -6.216216216216
$ws.Range("I16").Value2 = 1670
$ws.Range("J16").Value2 = 1688
$ws.Range("K16").Value2 = -1.0663507109
$ws.Range("L16").Value2 = 38.589211618257
$ws.Range("M16").Value2 = 5.629348513598
$ws.Range("N16").Value2 = -72.986088644451
# Row 17
$ws.Range("C17").Value2 = 144
$ws.Range("D17").Value2 = 166
$ws.Range("E17").Value2 = -13.253012048192
$ws.Range("F17").Value2 = 610
$ws.Range("G17").Value2 = 583
$ws.Range("H17").Value2 = 4.631217838765
$ws.Range("I17").Value2 = 2828
$ws.Range("J17").Value2 = 2569
$ws.Range("K17").Value2 = 10.081743869209
$ws.Range("L17").Value2 = 33.774834437086
$ws.Range("M17").Value2 = 72.12416311625
$ws.Range("N17").Value2 = -9.993634627625
# Row 18
$ws.Range("C18").Value2 = 58
$ws.Range("D18").Value2 = 57
$ws.Range("E18").Value2 = 1.754385964912
$ws.Range("F18").Value2 = 218
$ws.Range("G18").Value2 = 227
$ws.Range("H18").Value2 = -3.964757709251
$ws.Range("I18").Value2 = 1177
$ws.Range("J18").Value2 = 1150
$ws.Range("K18").Value2 = 2.347826086956
$ws.Range("L18").Value2 = 45.848822800495
$ws.Range("M18").Value2 = 2.170138888888
$ws.Range("N18").Value2 = -83.543064876957
# Row 19
$ws.Range("C19").Value2 = 159
$ws.Range("D19").Value2 = 158
$ws.Range("E19").Value2 = 0.632911392405
$ws.Range("F19").Value2 = 585
$ws.Range("G19").Value2 = 580
$ws.Range("H19").Value2 = 0.862068965517
$ws.Range("I19").Value2 = 2900
$ws.Range("J19").Value2 = 2997
$ws.Range("K19").Value2 = -3.236569903236
$ws.Range("L19").Value2 = 30.571814497973
$ws.Range("M19").Value2 = 78.242163491087
$ws.Range("N19").Value2 = 7.766629505759
# Row 20
$ws.Range("C20").Value2 = 91
$ws.Range("D20").Value2 = 74
$ws.Range("E20").Value2 = 22.972972972973
$ws.Range("F20").Value2 = 441
$ws.Range("G20").Value2 = 258
$ws.Range("H20").Value2 = 70.930232558139
$ws.Range("I20").Value2 = 2065
$ws.Range("J20").Value2 = 1558
$ws.Range("K20").Value2 = 32.541720154043
$ws.Range("L20").Value2 = 114.656964656965
$ws.Range("M20").Value2 = 172.427440633245
$ws.Range("N20").Value2 = -65.24739145069
# Row 21
$ws.Range("C21").Value2 = 563
$ws.Range("D21").Value2 = 551
$ws.Range("E21").Value2 = 2.177858439201
$ws.Range("F21").Value2 = 2246
$ws.Range("G21").Value2 = 2058
$ws.Range("H21").Value2 = 9.13508260447
$ws.Range("I21").Value2 = 10844
$ws.Range("J21").Value2 = 10171
$ws.Range("K21").Value2 = 6.616851833644
$ws.Range("L21").Value2 = 44.721740290938
$ws.Range("M21").Value2 = 56.931982633864
$ws.Range("N21").Value2 = -57.566034044218
# Row 22
$ws.Range("D22").Value2 = 14
$ws.Range("E22").Value2 = -42.857142857142
$ws.Range("G22").Value2 = 35
$ws.Range("H22").Value2 = -31.428571428571
$ws.Range("I22").Value2 = 113
$ws.Range("J22").Value2 = 138
$ws.Range("K22").Value2 = -18.115942028985
$ws.Range("L22").Value2 = 28.40909090909
$ws.Range("M22").Value2 = -12.403100775193
# Row 23
$ws.Range("C23").Value2 = 19
$ws.Range("D23").Value2 = 35
$ws.Range("E23").Value2 = -45.714285714285
$ws.Range("F23").Value2 = 121
$ws.Range("G23").Value2 = 127
$ws.Range("H23").Value2 = -4.724409448818
$ws.Range("I23").Value2 = 666
$ws.Range("J23").Value2 = 578
$ws.Range("K23").Value2 = 15.224913494809
$ws.Range("L23").Value2 = 54.524361948955
$ws.Range("M23").Value2 = 78.552278820375
# Row 24
$ws.Range("C24").Value2 = 355
$ws.Range("D24").Value2 = 365
$ws.Range("E24").Value2 = -2.739726027397
$ws.Range("F24").Value2 = 1371
$ws.Range("G24").Value2 = 1419
$ws.Range("H24").Value2 = -3.38266384778
$ws.Range("I24").Value2 = 6664
$ws.Range("J24").Value2 = 6734
$ws.Range("K24").Value2 = -1.039501039501
$ws.Range("L24").Value2 = 45.343511450381
$ws.Range("M24").Value2 = 45.311818578281
# Row 25
$ws.Range("C25").Value2 = 210
$ws.Range("D25").Value2 = 199
$ws.Range("E25").Value2 = 5.527638190954
$ws.Range("F25").Value2 = 880
$ws.Range("G25").Value2 = 821
$ws.Range("H25").Value2 = 7.186358099878
$ws.Range("I25").Value2 = 3927
$ws.Range("J25").Value2 = 3732
$ws.Range("K25").Value2 = 5.225080385852
$ws.Range("L25").Value2 = 32.088799192734
$ws.Range("M25").Value2 = -2.628316389784
# Row 26
$ws.Range("C26").Value2 = 8
$ws.Range("D26").Value2 = 14
$ws.Range("E26").Value2 = -42.857142857142
$ws.Range("F26").Value2 = 55
$ws.Range("G26").Value2 = 47
$ws.Range("H26").Value2 = 17.021276595744
$ws.Range("I26").Value2 = 251
$ws.Range("J26").Value2 = 265
$ws.Range("K26").Value2 = -5.283018867924
$ws.Range("L26").Value2 = 13.574660633484
# Row 27
$ws.Range("C27").Value2 = 25
$ws.Range("D27").Value2 = 23
$ws.Range("E27").Value2 = 8.695652173913
$ws.Range("F27").Value2 = 83
$ws.Range("G27").Value2 = 82
$ws.Range("H27").Value2 = 1.219512195121
$ws.Range("I27").Value2 = 404
$ws.Range("J27").Value2 = 348
$ws.Range("K27").Value2 = 16.091954022988
$ws.Range("L27").Value2 = 19.526627218934
# Row 28
$ws.Range("C28").Value2 = 8
$ws.Range("D28").Value2 = 13
$ws.Range("E28").Value2 = -38.461538461538
$ws.Range("F28").Value2 = 24
$ws.Range("G28").Value2 = 41
$ws.Range("H28").Value2 = -41.463414634146
$ws.Range("I28").Value2 = 123
$ws.Range("J28").Value2 = 195
$ws.Range("K28").Value2 = -36.923076923076
$ws.Range("L28").Value2 = -34.920634920634
$ws.Range("M28").Value2 = -19.078947368421
$ws.Range("N28").Value2 = -72.907488986784
# Row 29
$ws.Range("C29").Value2 = 8
$ws.Range("D29").Value2 = 11
$ws.Range("E29").Value2 = -27.272727272727
$ws.Range("F29").Value2 = 21
$ws.Range("G29").Value2 = 34
$ws.Range("H29").Value2 = -38.235294117647
$ws.Range("I29").Value2 = 100
$ws.Range("J29").Value2 = 166
$ws.Range("K29").Value2 = -39.759036144578
$ws.Range("L29").Value2 = -39.024390243902
$ws.Range("M29").Value2 = -21.875
$ws.Range("N29").Value2 = -75.308641975308
# Row 30
$ws.Range("C30").Value2 = 1
$ws.Range("E30").Value2 = -50
$ws.Range("F30").Value2 = 2
$ws.Range("G30").Value2 = 4
$ws.Range("I30").Value2 = 11
$ws.Range("J30").Value2 = 19
$ws.Range("K30").Value2 = -42.105263157894
$ws.Range("L30").Value2 = -59.259259259259

# Hate Crimes Week-to-Date 2023 value becomes numeric (was blank/"0" text) matching the D30 style
$ws.Range("C30").NumberFormat = $ws.Range("D30").NumberFormat
